# Agregar encabezados de la lista de articulos (fila 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "price"
$ws.Range("C1").Value = "stock"
$ws.Range("D1").Value = "code"
$ws.Range("E1").Value = "fechaActualizacion"
$ws.Range("F1").Value = "_id"

# Agregar el primer articulo (fila 2)
$ws.Range("A2").Value = "biopet perro adulto 20kg"

# precio y fecha "se ven" como numero/fecha pero deben quedar como texto
# (tal como los exporta la base de datos), por eso se escriben con un
# apostrofe inicial para forzar el tipo texto; luego se restaura el
# estilo "Normal" para que quede sin el formato de advertencia que Excel
# aplica automaticamente a numeros guardados como texto.
$ws.Range("B2").Value = "'23700"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = 44
$ws.Range("D2").Value = "ba1"

$ws.Range("E2").Value = "'10-09-2024"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "f557b6da-8ea9-465d-b92b-49921405c587"
